$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final grid of values for the A1:I4 range (header row + 3 data rows),
# matching the target sharedStrings/worksheet content described by the diff.
$values = @(
    @("Estado civil código", "Edad", "Estado civil", "Comarca nombre", "Número hogares", "Comarca código", "Provincia código", "Aragón", "Provincia nombre"),
    @("null", "iaest-measure:edad", "iaest-measure:numero-hogares", "null", "null", "sdmx-dimension:refArea", "null", "sdmx-dimension:refArea", "sdmx-dimension:refArea"),
    @("null", "medida", "medida", "dim", "medida", "null", "null", "dim", "dim"),
    @("null", "xsd:string", "xsd:string", "URI-comarca", "xsd:int", "null", "null", "URI-Comunidad", "URI-Provincia")
)

for ($r = 0; $r -lt $values.Length; $r++) {
    $row = $values[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $row[$c]
    }
}

$wb.Save()
